$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 208.75
$ws.Range("I5").Value = 70.166664
$ws.Range("J5").Value = 624.5
$ws.Range("K5").Value = 70.166664
$ws.Range("L5").Value = 624.5
$ws.Range("M5").Value = 44.833336
$ws.Range("N5").Value = -854.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1323.5555
$ws.Range("I28").Value = 917
$ws.Range("J28").Value = 2136.6667
$ws.Range("K28").Value = 917
$ws.Range("L28").Value = 2136.6667
$ws.Range("M28").Value = -432
$ws.Range("N28").Value = -3106.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1636.5103
$ws.Range("I132").Value = 1368.9487
$ws.Range("J132").Value = 2680
$ws.Range("K132").Value = 4106.8461
$ws.Range("L132").Value = 8040
$ws.Range("M132").Value = -1576.8461
$ws.Range("N132").Value = -13100

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1581.4117
$ws.Range("I138").Value = 1399
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 4197
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = 943
$ws.Range("N138").Value = -23780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2653.739
$ws.Range("I2").Value = 1529.2307
$ws.Range("J2").Value = 4115.6
$ws.Range("K2").Value = 1529.2307
$ws.Range("L2").Value = 4115.6
$ws.Range("M2").Value = -1416.2307
$ws.Range("N2").Value = -4341.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 500
$ws.Range("I110").Value = 500
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 500
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2653.739
$ws.Range("I116").Value = 1529.2307
$ws.Range("J116").Value = 4115.6
$ws.Range("K116").Value = 1529.2307
$ws.Range("L116").Value = 4115.6
$ws.Range("M116").Value = 764.7692999999999
$ws.Range("N116").Value = -8703.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3142.5557
$ws.Range("I122").Value = 2451
$ws.Range("J122").Value = 4525.6665
$ws.Range("K122").Value = 7353
$ws.Range("L122").Value = 13576.9995
$ws.Range("M122").Value = -4903
$ws.Range("N122").Value = -18476.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2653.739
$ws.Range("I3").Value = 1529.2307
$ws.Range("J3").Value = 4115.6
$ws.Range("K3").Value = 1529.2307
$ws.Range("L3").Value = 4115.6
$ws.Range("M3").Value = -1415.2307
$ws.Range("N3").Value = -4343.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2207.6155
$ws.Range("I134").Value = 2283.25
$ws.Range("J134").Value = 1300
$ws.Range("K134").Value = 6849.75
$ws.Range("L134").Value = 3900
$ws.Range("M134").Value = -4314.75
$ws.Range("N134").Value = -8970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 368.33334
$ws.Range("I16").Value = 368.33334
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 368.33334
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -81.33334000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1643.3334
$ws.Range("I31").Value = 1643.3334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1643.3334
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -1348.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1643.3334
$ws.Range("I34").Value = 1643.3334
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1643.3334
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -1441.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7527.9375
$ws.Range("I99").Value = 7563.1333
$ws.Range("J99").Value = 7000
$ws.Range("K99").Value = 7563.1333
$ws.Range("L99").Value = 7000
$ws.Range("M99").Value = -6065.1333
$ws.Range("N99").Value = -9996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 368.33334
$ws.Range("I113").Value = 368.33334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 368.33334
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1801.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7527.9375
$ws.Range("I126").Value = 7563.1333
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 22689.3999
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -20219.3999
$ws.Range("N126").Value = -25940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 7887.222
$ws.Range("I80").Value = 5996.75
$ws.Range("J80").Value = 9399.6
$ws.Range("K80").Value = 17990.25
$ws.Range("L80").Value = 28198.8
$ws.Range("M80").Value = -17054.25
$ws.Range("N80").Value = -30070.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 7887.222
$ws.Range("I83").Value = 5996.75
$ws.Range("J83").Value = 9399.6
$ws.Range("K83").Value = 53970.75
$ws.Range("L83").Value = 84596.40000000001
$ws.Range("M83").Value = -49290.75
$ws.Range("N83").Value = -93956.40000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 39841.69
$ws.Range("I112").Value = 19316.666
$ws.Range("J112").Value = 45999.2
$ws.Range("K112").Value = 57949.99800000001
$ws.Range("L112").Value = 137997.6
$ws.Range("M112").Value = -56841.99800000001
$ws.Range("N112").Value = -140213.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 853.6667
$ws.Range("I113").Value = 774
$ws.Range("J113").Value = 893.5
$ws.Range("K113").Value = 2322
$ws.Range("L113").Value = 2680.5
$ws.Range("M113").Value = -152
$ws.Range("N113").Value = -7020.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4021.3
$ws.Range("I122").Value = 4031.389
$ws.Range("J122").Value = 3930.5
$ws.Range("K122").Value = 12094.167
$ws.Range("L122").Value = 11791.5
$ws.Range("M122").Value = -9644.167000000001
$ws.Range("N122").Value = -16691.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1769.4694
$ws.Range("I132").Value = 1590.7446
$ws.Range("J132").Value = 5969.5
$ws.Range("K132").Value = 4772.2338
$ws.Range("L132").Value = 17908.5
$ws.Range("M132").Value = -2242.2338
$ws.Range("N132").Value = -22968.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4674
$ws.Range("I40").Value = 4674
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4674
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1492.8125
$ws.Range("I132").Value = 1468.1538
$ws.Range("J132").Value = 1599.6666
$ws.Range("K132").Value = 4404.4614
$ws.Range("L132").Value = 4798.9998
$ws.Range("M132").Value = -1874.4614
$ws.Range("N132").Value = -9858.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 11662.75
$ws.Range("I55").Value = 2049.5
$ws.Range("J55").Value = 21276
$ws.Range("K55").Value = 2049.5
$ws.Range("L55").Value = 21276
$ws.Range("M55").Value = -1772.5
$ws.Range("N55").Value = -21830

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 662.5454999999999
$ws.Range("I107").Value = 624.875
$ws.Range("J107").Value = 763
$ws.Range("K107").Value = 1874.625
$ws.Range("L107").Value = 2289
$ws.Range("M107").Value = 45.375
$ws.Range("N107").Value = -6129

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3800
$ws.Range("I122").Value = 3800
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11400
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1929.0625
$ws.Range("I136").Value = 1993.5714
$ws.Range("J136").Value = 1477.5
$ws.Range("K136").Value = 5980.7142
$ws.Range("L136").Value = 4432.5
$ws.Range("M136").Value = -3430.7142
$ws.Range("N136").Value = -9532.5
